# Actualización automática 2025-06-19 14:50:09
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Narrow column D slightly (stored OOXML width 13 -> 12)
$ws.Columns.Item(4).ColumnWidth = 11.1666666666667

# Row 2 (OTROS)
$ws.Range("D2").Value = 3566
$ws.Range("E2").Value = -3566

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 321.85
$ws.Range("E3").Value = 17178.15
$ws.Range("F3").Value = 0.01839142857142857

# Row 4 (TOTAL)
$ws.Range("D4").Value = 3887.85
$ws.Range("E4").Value = 13612.15
$ws.Range("F4").Value = 0.2221628571428571
